$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 32
$ws.Cells.Item(32, 8).Value = 4978807
$ws.Cells.Item(32, 9).Value = 349.5
$ws.Cells.Item(32, 11).Value = 349.5
$ws.Cells.Item(32, 13).Value = -23.5

# ALC row 137
$ws.Cells.Item(137, 8).Value = 401799.22
$ws.Cells.Item(137, 9).Value = 3580.0476
$ws.Cells.Item(137, 10).Value = 1098682.8
$ws.Cells.Item(137, 11).Value = 10740.1428
$ws.Cells.Item(137, 12).Value = 3296048.4
$ws.Cells.Item(137, 13).Value = -8190.1428
$ws.Cells.Item(137, 14).Value = -3301148.4

# ALC row 138
$ws.Cells.Item(138, 8).Value = 4552.2334
$ws.Cells.Item(138, 9).Value = 3974.625
$ws.Cells.Item(138, 10).Value = 4641.096
$ws.Cells.Item(138, 11).Value = 11923.875
$ws.Cells.Item(138, 12).Value = 13923.288
$ws.Cells.Item(138, 13).Value = -6783.875
$ws.Cells.Item(138, 14).Value = -24203.288

$ws = $wb.Worksheets.Item("ARM")
# ARM row 74
$ws.Cells.Item(74, 8).Value = 10932327
$ws.Cells.Item(74, 9).Value = 7804434.5
$ws.Cells.Item(74, 10).Value = 18578286
$ws.Cells.Item(74, 11).Value = 7804434.5
$ws.Cells.Item(74, 12).Value = 18578286
$ws.Cells.Item(74, 13).Value = -7803560.5
$ws.Cells.Item(74, 14).Value = -18580034

# ARM row 77
$ws.Cells.Item(77, 8).Value = 10932327
$ws.Cells.Item(77, 9).Value = 7804434.5
$ws.Cells.Item(77, 10).Value = 18578286
$ws.Cells.Item(77, 11).Value = 39022172.5
$ws.Cells.Item(77, 12).Value = 92891430
$ws.Cells.Item(77, 13).Value = -39017804.5
$ws.Cells.Item(77, 14).Value = -92900166

# ARM row 110
$ws.Cells.Item(110, 8).Value = 549.86664
$ws.Cells.Item(110, 9).Value = 549.86664
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 549.86664
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = 1495.13336
$ws.Cells.Item(110, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 82
$ws.Cells.Item(82, 8).Value = 17360
$ws.Cells.Item(82, 9).Value = 6776
$ws.Cells.Item(82, 10).Value = 35000
$ws.Cells.Item(82, 11).Value = 6776
$ws.Cells.Item(82, 12).Value = 35000
$ws.Cells.Item(82, 13).Value = -6393
$ws.Cells.Item(82, 14).Value = -35766

# BSM row 85
$ws.Cells.Item(85, 8).Value = 17360
$ws.Cells.Item(85, 9).Value = 6776
$ws.Cells.Item(85, 10).Value = 35000
$ws.Cells.Item(85, 11).Value = 6776
$ws.Cells.Item(85, 12).Value = 35000
$ws.Cells.Item(85, 13).Value = -5450
$ws.Cells.Item(85, 14).Value = -37652

# BSM row 107
$ws.Cells.Item(107, 8).Value = 1462.7894
$ws.Cells.Item(107, 9).Value = 1494.1
$ws.Cells.Item(107, 10).Value = 1345.375
$ws.Cells.Item(107, 11).Value = 1494.1
$ws.Cells.Item(107, 12).Value = 1345.375
$ws.Cells.Item(107, 13).Value = 425.9000000000001
$ws.Cells.Item(107, 14).Value = -5185.375

# BSM row 134
$ws.Cells.Item(134, 8).Value = 35147.094
$ws.Cells.Item(134, 9).Value = 57028.723
$ws.Cells.Item(134, 10).Value = 7013.5713
$ws.Cells.Item(134, 11).Value = 171086.169
$ws.Cells.Item(134, 12).Value = 21040.7139
$ws.Cells.Item(134, 13).Value = -168551.169
$ws.Cells.Item(134, 14).Value = -26110.7139

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Cells.Item(22, 8).Value = 1079.875
$ws.Cells.Item(22, 9).Value = 781.6667
$ws.Cells.Item(22, 10).Value = 1974.5
$ws.Cells.Item(22, 11).Value = 781.6667
$ws.Cells.Item(22, 12).Value = 1974.5
$ws.Cells.Item(22, 13).Value = -431.6667
$ws.Cells.Item(22, 14).Value = -2674.5

# CRP row 31
$ws.Cells.Item(31, 8).Value = 2913.4727
$ws.Cells.Item(31, 9).Value = 2199.3333
$ws.Cells.Item(31, 10).Value = 4654.1875
$ws.Cells.Item(31, 11).Value = 2199.3333
$ws.Cells.Item(31, 12).Value = 4654.1875
$ws.Cells.Item(31, 13).Value = -1904.3333
$ws.Cells.Item(31, 14).Value = -5244.1875

# CRP row 34
$ws.Cells.Item(34, 8).Value = 2913.4727
$ws.Cells.Item(34, 9).Value = 2199.3333
$ws.Cells.Item(34, 10).Value = 4654.1875
$ws.Cells.Item(34, 11).Value = 2199.3333
$ws.Cells.Item(34, 12).Value = 4654.1875
$ws.Cells.Item(34, 13).Value = -1997.3333
$ws.Cells.Item(34, 14).Value = -5058.1875

# CRP row 141
$ws.Cells.Item(141, 8).Value = 29432.691
$ws.Cells.Item(141, 9).Value = 31500
$ws.Cells.Item(141, 10).Value = 28513.889
$ws.Cells.Item(141, 11).Value = 31500
$ws.Cells.Item(141, 12).Value = 28513.889
$ws.Cells.Item(141, 13).Value = -26320
$ws.Cells.Item(141, 14).Value = -38873.889

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Cells.Item(132, 8).Value = 4452.8335
$ws.Cells.Item(132, 9).Value = 5071.0625
$ws.Cells.Item(132, 10).Value = 3746.2856
$ws.Cells.Item(132, 11).Value = 15213.1875
$ws.Cells.Item(132, 12).Value = 11238.8568
$ws.Cells.Item(132, 13).Value = -12683.1875
$ws.Cells.Item(132, 14).Value = -16298.8568

$ws = $wb.Worksheets.Item("LTW")
# LTW row 61
$ws.Cells.Item(61, 8).Value = 1690.7826
$ws.Cells.Item(61, 9).Value = 1688.0588
$ws.Cells.Item(61, 10).Value = 1698.5
$ws.Cells.Item(61, 11).Value = 1688.0588
$ws.Cells.Item(61, 12).Value = 1698.5
$ws.Cells.Item(61, 13).Value = -1486.0588
$ws.Cells.Item(61, 14).Value = -2102.5

# LTW row 113
$ws.Cells.Item(113, 8).Value = 1690.7826
$ws.Cells.Item(113, 9).Value = 1688.0588
$ws.Cells.Item(113, 10).Value = 1698.5
$ws.Cells.Item(113, 11).Value = 1688.0588
$ws.Cells.Item(113, 12).Value = 1698.5
$ws.Cells.Item(113, 13).Value = 481.9412
$ws.Cells.Item(113, 14).Value = -6038.5

# LTW row 132
$ws.Cells.Item(132, 8).Value = 6295805
$ws.Cells.Item(132, 9).Value = 2107.6667
$ws.Cells.Item(132, 10).Value = 14504975
$ws.Cells.Item(132, 11).Value = 6323.000100000001
$ws.Cells.Item(132, 12).Value = 43514925
$ws.Cells.Item(132, 13).Value = -3793.000100000001
$ws.Cells.Item(132, 14).Value = -43519985

# LTW row 133
$ws.Cells.Item(133, 8).Value = 38933.332
$ws.Cells.Item(133, 10).Value = 38933.332
$ws.Cells.Item(133, 12).Value = 38933.332
$ws.Cells.Item(133, 14).Value = -43993.332

$ws = $wb.Worksheets.Item("WVR")
# WVR row 4
$ws.Cells.Item(4, 8).Value = 6491.2856
$ws.Cells.Item(4, 9).Value = 4500
$ws.Cells.Item(4, 10).Value = 6823.1665
$ws.Cells.Item(4, 11).Value = 4500
$ws.Cells.Item(4, 12).Value = 6823.1665
$ws.Cells.Item(4, 13).Value = -4387
$ws.Cells.Item(4, 14).Value = -7049.1665

# WVR row 62
$ws.Cells.Item(62, 8).Value = 4565.3335
$ws.Cells.Item(62, 9).Value = 4990
$ws.Cells.Item(62, 11).Value = 4990
$ws.Cells.Item(62, 13).Value = -4366

# WVR row 65
$ws.Cells.Item(65, 8).Value = 4565.3335
$ws.Cells.Item(65, 9).Value = 4990
$ws.Cells.Item(65, 11).Value = 24950
$ws.Cells.Item(65, 13).Value = -21830

# WVR row 93
$ws.Cells.Item(93, 8).Value = 30000
$ws.Cells.Item(93, 10).Value = 30000
$ws.Cells.Item(93, 12).Value = 30000
$ws.Cells.Item(93, 14).Value = -34992

# WVR row 96
$ws.Cells.Item(96, 8).Value = 22999
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 22999
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 22999
$ws.Cells.Item(96, 13).ClearContents()
$ws.Cells.Item(96, 14).Value = -25745

# WVR row 132
$ws.Cells.Item(132, 8).Value = 2929.2285
$ws.Cells.Item(132, 9).Value = 1975.2778
$ws.Cells.Item(132, 10).Value = 3939.2942
$ws.Cells.Item(132, 11).Value = 5925.8334
$ws.Cells.Item(132, 12).Value = 11817.8826
$ws.Cells.Item(132, 13).Value = -3395.8334
$ws.Cells.Item(132, 14).Value = -16877.8826

# WVR row 136
$ws.Cells.Item(136, 8).Value = 263992.88
$ws.Cells.Item(136, 9).Value = 38500.223
$ws.Cells.Item(136, 10).Value = 669879.7
$ws.Cells.Item(136, 11).Value = 115500.669
$ws.Cells.Item(136, 12).Value = 2009639.1
$ws.Cells.Item(136, 13).Value = -112950.669
$ws.Cells.Item(136, 14).Value = -2014739.1
